# fdo#42624: add range names referencing range names that are loaded
# after the referencing range name.
#
# Adds two new global defined names:
#   Global5 -> Global6            (forward reference to a name defined later)
#   Global6 -> Sheet2!$B$1
# and a new formula cell on Sheet2 (A6) that exercises the forward
# reference via Global5.

$wb = $excel.ActiveWorkbook

# New global (workbook-scope) defined names, inserted right after Global4.
$wb.Names.Add("Global5", "=Global6") | Out-Null
$wb.Names.Add("Global6", "=Sheet2!`$B`$1") | Out-Null

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New formula cell referencing the forward-declared name.
$ws2.Range("A6").Formula = "=Global5"

# Restore/update the cell selections on both sheets to match the saved
# view state (Sheet1 collapses the old A4:C5 block selection down to a
# single cell, Sheet2's selection moves on to the newly added row).
$ws1.Range("A5").Select() | Out-Null
$ws2.Range("A7").Select() | Out-Null
